$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Coinranking snapshot refresh: update prices / 1h volume deltas for every
# row, plus the HuobiToken <-> TrustWalletToken row swap at 43/45.
# A leading "'" keeps Price (column D) values that look like plain numbers
# (e.g. "231.56") stored as text, matching the source data (prices such as
# "37.430.78" use dotted thousands separators and already stay text as-is).
$ws.Range("D2").Value = "37.430.78"
$ws.Range("E2").Value = "  -1.25%  "
$ws.Range("D3").Value = "2.057.87"
$ws.Range("E3").Value = "  -1.31%  "
$ws.Range("E4").Value = "  -0.12%  "
$ws.Range("D5").Value = "'231.56"
$ws.Range("E5").Value = "  -0.76%  "
$ws.Range("E6").Value = "  -0.61%  "
$ws.Range("D8").Value = "'57.30"
$ws.Range("E8").Value = "  -3.46%  "
$ws.Range("E9").Value = "  -2.84%  "
$ws.Range("D10").Value = "'0.0773"
$ws.Range("E10").Value = "  -2.02%  "
$ws.Range("E11").Value = "  +1.40%  "
$ws.Range("D12").Value = "2.360.03"
$ws.Range("E12").Value = "  -1.37%  "
$ws.Range("D13").Value = "'14.64"
$ws.Range("E13").Value = "  -0.62%  "
$ws.Range("D14").Value = "'21.08"
$ws.Range("E14").Value = "  -0.67%  "
$ws.Range("D15").Value = "'0.760"
$ws.Range("E15").Value = "  -2.16%  "
$ws.Range("D16").Value = "'5.32"
$ws.Range("E16").Value = "  -0.69%  "
$ws.Range("D17").Value = "2.056.72"
$ws.Range("E17").Value = "  -2.29%  "
$ws.Range("D18").Value = "37.424.55"
$ws.Range("E18").Value = "  -1.05%  "
$ws.Range("E19").Value = "  -0.79%  "
$ws.Range("D20").Value = "'69.74"
$ws.Range("E20").Value = "  -2.72%  "
$ws.Range("D21").Value = "0.0₃0824"
$ws.Range("E21").Value = "  -3.06%  "
$ws.Range("D22").Value = "'226.65"
$ws.Range("E22").Value = "  -0.71%  "
$ws.Range("E24").Value = "  +0.37%  "
$ws.Range("E25").Value = "  -2.90%  "
$ws.Range("D26").Value = "'9.88"
$ws.Range("E26").Value = "  +7.54%  "
$ws.Range("D27").Value = "'170.26"
$ws.Range("E27").Value = "  -1.09%  "
$ws.Range("E28").Value = "  -5.88%  "
$ws.Range("D29").Value = "'19.22"
$ws.Range("E29").Value = "  -1.48%  "
$ws.Range("E30").Value = "  -5.22%  "
$ws.Range("E31").Value = "  +0.15%  "
$ws.Range("D32").Value = "'4.54"
$ws.Range("E32").Value = "  -4.06%  "
$ws.Range("E33").Value = "  -1.56%  "
$ws.Range("D34").Value = "'4.59"
$ws.Range("E34").Value = "  -2.91%  "
$ws.Range("D35").Value = "'2.50"
$ws.Range("E35").Value = "  -0.44%  "
$ws.Range("D36").Value = "'1.83"
$ws.Range("E36").Value = "  +0.17%  "
$ws.Range("E37").Value = "  -4.32%  "
$ws.Range("E38").Value = "  -0.12%  "
$ws.Range("E39").Value = "  -1.86%  "
$ws.Range("D40").Value = "'0.0225"
$ws.Range("E40").Value = "  +3.04%  "
$ws.Range("D41").Value = "'98.48"
$ws.Range("E41").Value = "  -0.68%  "
$ws.Range("D42").Value = "'0.0958"
$ws.Range("E42").Value = "  -2.61%  "
$ws.Range("B43").Value = "TrustWalletToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D43").Value = "'1.19"
$ws.Range("E43").Value = "  +3.36%  "
$ws.Range("D44").Value = "1.477.53"
$ws.Range("E44").Value = "  +2.03%  "
$ws.Range("B45").Value = "HuobiToken"
$ws.Range("C45").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D45").Value = "'2.89"
$ws.Range("E45").Value = "  -0.41%  "
$ws.Range("D46").Value = "'16.62"
$ws.Range("E46").Value = "  -1.33%  "
$ws.Range("E47").Value = "  -2.92%  "
$ws.Range("D48").Value = "'7.25"
$ws.Range("E48").Value = "  -1.84%  "
$ws.Range("E49").Value = "  -6.24%  "
$ws.Range("D50").Value = "'2.96"
$ws.Range("E50").Value = "  -1.41%  "
$ws.Range("D51").Value = "2.246.57"
$ws.Range("E51").Value = "  -1.35%  "
